$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: student milestone changes from "I" to "II"
$ws.Range("E5").Value = "II"

# Row 15: clear the "Milestone Completed" (X) flag in column F
$ws.Range("F15").ClearContents()

# Row 41: mark student milestone as "III"
$ws.Range("E41").Value = "III"

# Row 63: change student milestone from "II" to "III"
$ws.Range("E63").Value = "III"

# Row 67: mark student milestone as "III"
$ws.Range("E67").Value = "III"

# Row 85: Island (Tropical) Theme row - mark student milestone as "III"
$ws.Range("E85").Value = "III"

# Row 97: add new project source citation (coconut tree / island asset used for the Island theme)
$ws.Range("A97").Value = "https://www.cgtrader.com/free-3d-models/exterior/landscape/coconut-tree-island-low-poly"

# Update the visible scroll position / selection to reflect where the user ended up working
$ws.Application.ActiveWindow.ScrollRow = 88
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A97").Select()
